$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column G readings (rows 3-26) ---
$gValues = @{
    3  = 53
    4  = 51
    5  = 51
    6  = 51
    7  = 51
    8  = 51
    9  = 51
    10 = 51
    11 = 51
    12 = 51
    13 = 51
    14 = 51
    15 = 51
    16 = 51
    17 = 51
    18 = 51
    19 = 55
    20 = 51
    21 = 51
    22 = 51
    23 = 51
    24 = 51
    25 = 51
    26 = 51
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# --- Update the sheet view (scroll position + selection) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("G20").Select()
